$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction in SA algorithm: recompute the "Fitness" (best-so-far) column
# for this run's log. Values below reflect the corrected staircase of
# best fitness achieved per generation.
$ws.Range("C2:C8").Value = 8127
$ws.Range("C9:C26").Value = 7773
$ws.Range("C27:C50").Value = 7318
$ws.Range("C51:C109").Value = 7310
